# Insert a new weekly data row at row 130 (Chirimoya, Macroferia Regional de Talca),
# pushing all subsequent rows (old 130..162) down by one (new 131..163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert a blank row above the current row 130.
$ws.Rows.Item(130).Insert()

# Fill in the new row's values.
$ws.Range("A130").Value = 5
$ws.Range("B130").Value = "Macroferia Regional de Talca"
$ws.Range("C130").Value = "Maule"
$ws.Range("D130").Value = 45204
$ws.Range("E130").Value = 7
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100107
$ws.Range("H130").Value = "Otros"
$ws.Range("I130").Value = 100107002
$ws.Range("J130").Value = "Chirimoya"
$ws.Range("K130").Value = "Cultivar IV Región"
$ws.Range("L130").Value = "Segunda"
$ws.Range("M130").Value = 150
$ws.Range("N130").Value = 20000
$ws.Range("O130").Value = 20000
$ws.Range("P130").Value = 20000
$ws.Range("Q130").Value = "$/bandeja 10 kilos"
$ws.Range("R130").Value = "Provincia de Limarí"
$ws.Range("S130").Value = 2000
$ws.Range("T130").Value = 10
